# Rename ConsoleCatchall => Reconciliate
# Underlying data tweaks that accompanied the rename in the source workbook.

$wb = $excel.ActiveWorkbook

# --- TestRecord sheet: bump the sample transaction date & amount by one day / 1.20 ---
$wsTestRecord = $wb.Worksheets.Item("TestRecord")
$wsTestRecord.Range("A10").Value = 43218
$wsTestRecord.Range("B10").Value = 67.14
$wsTestRecord.Range("E10").Value = "some test textzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Budget Out sheet: corresponding amount bump ---
$wsBudgetOut = $wb.Worksheets.Item("Budget Out")
$wsBudgetOut.Range("C9").Value = 83.02
$wsBudgetOut.Range("F9").Value = "Description007zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"

# --- Expected Out sheet: matching totals update (B1 SUM formula recalculates automatically) ---
$wsExpectedOut = $wb.Worksheets.Item("Expected Out")
$wsExpectedOut.Range("B9").Value = 1324.16
$wsExpectedOut.Range("B11").Value = 420.82
